# Updates cryptos list values (price/volume) per upstream data refresh.
# D-column "Price" cells are plain literal text (e.g. "22.388.07"), not
# numbers -- Excel auto-parses bare numeric-looking strings into real
# numbers/dates, which would corrupt both the displayed text (dropping
# trailing zeros, adding float noise) and the cell type. Prefixing with
# a leading apostrophe forces literal text entry (Excel's own quote-
# prefix convention); resetting .Style to "Normal" afterwards drops the
# transient quote-prefix style flag so the cell format matches the rest
# of the sheet (no explicit style).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'22.390.62"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.31%  "
$ws.Range("D3").Value = "'1.561.79"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.74%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("E5").Value = "  -0.13%  "
$ws.Range("D6").Value = "'286.04"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.77%  "
$ws.Range("D7").Value = "'0.3642"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.65%  "
$ws.Range("D8").Value = "'48.44"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.09%  "
$ws.Range("D9").Value = "'0.3337"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.08%  "
$ws.Range("D10").Value = "'1.126"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.43%  "
$ws.Range("D11").Value = "'0.07405"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.21%  "
$ws.Range("E12").Value = "  -0.13%  "
$ws.Range("D13").Value = "'20.79"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.65%  "
$ws.Range("D14").Value = "'5.920"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.20%  "
$ws.Range("D15").Value = "'6.876"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.99%  "
$ws.Range("D16").Value = "'1.561.29"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.00%  "
$ws.Range("D17").Value = "'0.00001102"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.00%  "
$ws.Range("D18").Value = "'88.73"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.49%  "
$ws.Range("D19").Value = "'0.06726"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.10%  "
$ws.Range("E20").Value = "  -0.14%  "
$ws.Range("D21").Value = "'6.324"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.81%  "
$ws.Range("D22").Value = "'16.04"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.45%  "
$ws.Range("D23").Value = "'11.93"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.24%  "
$ws.Range("D24").Value = "'22.376.04"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.39%  "
$ws.Range("D25").Value = "'2.401"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.73%  "
$ws.Range("D26").Value = "'2.537"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.87%  "
$ws.Range("D27").Value = "'149.29"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.60%  "
$ws.Range("D28").Value = "'19.39"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.71%  "
$ws.Range("D29").Value = "'4.996"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.31%  "
$ws.Range("D30").Value = "'122.87"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.44%  "
$ws.Range("D31").Value = "'1.734.74"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.89%  "
$ws.Range("D32").Value = "'1.059"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.17%  "
$ws.Range("D33").Value = "'6.114"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.22%  "
$ws.Range("D34").Value = "'1.979"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("D35").Value = "'9.581"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.71%  "
$ws.Range("D36").Value = "'0.08229"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.03%  "
$ws.Range("D37").Value = "'0.02382"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.15%  "
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "'0.06361"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.51%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").Value = "'1.295"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.32%  "
$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D40").Value = "'0.2208"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.91%  "
$ws.Range("D41").Value = "'5.327"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.83%  "
$ws.Range("D42").Value = "'11.12"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.12%  "
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").Value = "'0.6057"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.51%  "
$ws.Range("B44").Value = "Frax"
$ws.Range("C44").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D44").Value = "'1.000"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.12%  "
$ws.Range("E45").Value = "  -2.45%  "
$ws.Range("D46").Value = "'3.757"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.46%  "
$ws.Range("D47").Value = "'0.5736"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.07%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "'2.002"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.15%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "'124.37"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.03%  "
$ws.Range("D50").Value = "'1.208"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.86%  "
$ws.Range("D51").Value = "'0.07217"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.55%  "
